$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$notes = $s.NotesPage
$notes.Shapes.AddPlaceholder(2)
$notes.HeadersFooters.SlideNumber.Visible = $true
